$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.986.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.769.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3521"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.92"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07397"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.085"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.018"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.202"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.766.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06422"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.93"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.801"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.030.58"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.15"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.158"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.66"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.970.03"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.30"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.077"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.656"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.564"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.70"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02271"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06120"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2072"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.914"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.194"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6165"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.444"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.779"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.747"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5810"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.03"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.936"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.128"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06811"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.12"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.54%  "
